$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 16) duplicating the pattern of row 15, but for the new
# "toy-next-question-id" entry (nb = 14), to be compatible with topic management.
$ws.Range("A16").Value = "14"
$ws.Range("B16").Value = "toy-next-question-id"
$ws.Range("C16").Value = "Curious"
$ws.Range("D16").Value = "happy"
$ws.Range("E16").Value = 0.75
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 80
$ws.Range("H16").Value = "MAIN-SIT_1-R_Likesth"
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 80
$ws.Range("K16").Value = 1

# Match the styles used by the rest of the column (copy formats from the row above,
# which already has the correct text/percent/general number formats per column).
$ws.Range("A15:K15").Copy()
$ws.Range("A16:K16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Keep the same explicit row height as the other data rows.
$ws.Rows.Item(16).RowHeight = $ws.Rows.Item(15).RowHeight

# Update the active selection to the newly added cell, like the author left it.
$ws.Range("B16").Select()
